$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.800.85'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = '3.575.71'
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '189.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("E7").Value = '  -3.15%  '
$ws.Range("D8").Value = '3.571.64'
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -1.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.662'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.76'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000302'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").Value = '4.149.66'
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.79'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.96%  '
$ws.Range("D17").Value = '3.564.17'
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").Value = '69.742.86'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '475.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +11.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.17%  '
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '95.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.91%  '
$ws.Range("E27").Value = '  -3.98%  '
$ws.Range("E28").Value = '  -2.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.119'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.06'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '583.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.24%  '
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.396'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.22'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +17.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.138'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.93%  '
$ws.Range("D43").Value = '3.236.25'
$ws.Range("E43").Value = '  -2.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.04%  '
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("E46").Value = '  -1.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.46%  '
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.96%  '
